$d = $word.ActiveDocument

$d.Content.Find.Execute("293÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "123÷6=", 2)
$d.Content.Find.Execute("133÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "995÷3=", 2)
$d.Content.Find.Execute("566÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "320÷3=", 2)
$d.Content.Find.Execute("811÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "131÷7=", 2)
$d.Content.Find.Execute("295÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "570÷8=", 2)
$d.Content.Find.Execute("361÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "356÷9=", 2)
$d.Content.Find.Execute("102÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "385÷4=", 2)
$d.Content.Find.Execute("269÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "758÷6=", 2)
$d.Content.Find.Execute("266÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "623÷6=", 2)
$d.Content.Find.Execute("773÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "743÷2=", 2)
$d.Content.Find.Execute("626÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "555÷5=", 2)
$d.Content.Find.Execute("355÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "513÷3=", 2)
$d.Content.Find.Execute("138÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "305÷4=", 2)
$d.Content.Find.Execute("112÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "495÷8=", 2)
$d.Content.Find.Execute("985÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "367÷4=", 2)
$d.Content.Find.Execute("802÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "477÷9=", 2)
$d.Content.Find.Execute("661÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "578÷3=", 2)
$d.Content.Find.Execute("548÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "928÷4=", 2)
$d.Content.Find.Execute("446÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "551÷8=", 2)
$d.Content.Find.Execute("659÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "382÷5=", 2)
$d.Content.Find.Execute("369÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "112÷3=", 2)
$d.Content.Find.Execute("849÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "834÷6=", 2)
$d.Content.Find.Execute("108÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "608÷3=", 2)
$d.Content.Find.Execute("687÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "604÷7=", 2)
$d.Content.Find.Execute("418÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "689÷4=", 2)
